# LH: Added 2 papers
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Main")

# ---- Row 24: NusG single-molecule tracking paper ----
$cell = $ws.Range("B24")
$cell.Value = "Single-molecule tracking reveals the functional allocation, in vivo interactions, and spatial organization of universal transcription factor NusG"
$cell.VerticalAlignment = -4108
$cell.WrapText = $true
$chars = $cell.Characters(61, 7)
$chars.Font.Italic = $true
$chars.Font.ColorIndex = -4105

$ws.Range("C24").Value = 2024
$ws.Range("D24").Value = "https://doi.org/10.1016/j.molcel.2024.01.025"
$ws.Hyperlinks.Add($ws.Range("D24"), "https://doi.org/10.1016/j.molcel.2024.01.025") | Out-Null

$ws.Range("F24").Value = 5
$ws.Range("G24").Value = 5
$ws.Range("H24").Value = 5
$ws.Range("I24").Value = 4
$ws.Range("J24").Value = 3
$ws.Range("K24").Value = "?"
$ws.Range("L24").Value = "?"
$ws.Range("M24").Value = "I left evaluating quality of experimental desing until the end because I felt that the points I was going to make in it fit better in the next three categories. When I returned to it I mostly reiterated the points I had made previously. I am also not sure how we are deciding the values to put int the replicability and reproducibility category - is it an average of the blue columns or just a general observation?"

# ---- Row 25: single molecule microscopy repair foci paper ----
$ws.Range("B25").Value = "Single molecule microscopy reveals key physical features of repair foci in living cells "
$ws.Range("C25").Value = 2021
$ws.Range("D25").Value = "https://doi.org/10.7554/eLife.60577"
$ws.Hyperlinks.Add($ws.Range("D25"), "https://doi.org/10.7554/eLife.60577") | Out-Null

$ws.Range("F25").Value = 5
$ws.Range("G25").Value = 4
$ws.Range("H25").Value = 4
$ws.Range("I25").Value = 5
$ws.Range("J25").Value = 1
$ws.Range("K25").Value = "?"
$ws.Range("L25").Value = "?"

# ---- Row 29: Cancer cell membrane-coated nanoparticles ----
$ws.Range("B29").Value = "Cancer Cell Membrane-Coated Nanoparticles for Anticancer Vaccination and Drug Delivery"
$ws.Range("C29").Value = 2014
$ws.Range("D29").Value = "https://pubs.acs.org/doi/full/10.1021/nl500618u"
$ws.Hyperlinks.Add($ws.Range("D29"), "https://pubs.acs.org/doi/full/10.1021/nl500618u") | Out-Null
$ws.Range("E29").Value = "Ioan Duchastel"

$ws.Range("F29").Value = 2
$ws.Range("G29").Value = 5
$ws.Range("H29").Value = 5
$ws.Range("I29").Value = 4
$ws.Range("J29").Value = 4
$ws.Range("K29").Value = "N/A (?)"
$ws.Range("L29").Value = 4
$ws.Range("M29").Value = "Qualitative study with a lot of details on the experimental procedure, to the reagents used in every steps to the tools used in the experiment. Not 5/5 as some of the tools used are expensive but not a big deal. Image analysis can be done on open source software so overall pretty thorough and good study to follow"
$ws.Rows.Item(29).RowHeight = 21.6

# ---- Row 30: Effective cancer targeting macrophage membrane ----
$ws.Range("B30").Value = "Effective cancer targeting and imaging using macrophage membranecamouflaged upconversion nanoparticles"
$ws.Range("C30").Value = 2016
$ws.Range("D30").Value = "https://onlinelibrary.wiley.com/doi/abs/10.1002/jbm.a.35927"
$ws.Hyperlinks.Add($ws.Range("D30"), "https://onlinelibrary.wiley.com/doi/abs/10.1002/jbm.a.35927") | Out-Null
$ws.Range("E30").Value = "Ioan Duchastel"

$ws.Range("F30").Value = 2
$ws.Range("G30").Value = 4
$ws.Range("H30").Value = 3
$ws.Range("I30").Value = 5
$ws.Range("J30").Value = 4
$ws.Range("K30").Value = "N/A (?)"
$ws.Range("L30").Value = 4
$ws.Range("M30").Value = "Qualitative study with a decent amount of detail, methods are just enough for this paper, you can still access them through other related papers (which seems to be to get more citations than anything) but reagent concentrations are severely lacking, making it quite difficult to reproduce."
$ws.Rows.Item(30).RowHeight = 27

# ---- Rows 31-33: student name "Ioan" ----
$ws.Range("E31").Value = "Ioan"
$ws.Range("E32").Value = "Ioan"
$ws.Range("E33").Value = "Ioan"

# ---- View state ----
$ws.Application.ActiveWindow.ScrollRow = 14
$ws.Application.ActiveWindow.ScrollColumn = 8
$ws.Range("M25").Select() | Out-Null
